$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.276.07"
$ws.Range("E2").Value = "  +3.73%  "

$ws.Range("D3").Value = "2.641.37"
$ws.Range("E3").Value = "  +2.85%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'598.66"
$ws.Range("E5").Value = "  +2.36%  "

$ws.Range("D6").Value = "'156.03"
$ws.Range("E6").Value = "  +5.62%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +1.98%  "

$ws.Range("D9").Value = "'0.119"
$ws.Range("E9").Value = "  +10.71%  "

$ws.Range("D10").Value = "'0.405"
$ws.Range("E10").Value = "  +6.76%  "

$ws.Range("D11").Value = "'5.82"
$ws.Range("E11").Value = "  +2.61%  "

$ws.Range("D12").Value = "'0.154"
$ws.Range("E12").Value = "  +2.71%  "

$ws.Range("D13").Value = "'29.26"
$ws.Range("E13").Value = "  +7.81%  "

$ws.Range("D14").Value = "'0.0000193"
$ws.Range("E14").Value = "  +28.03%  "

$ws.Range("D15").Value = "3.116.45"
$ws.Range("E15").Value = "  +2.88%  "

$ws.Range("D16").Value = "65.135.35"
$ws.Range("E16").Value = "  +3.79%  "

$ws.Range("D17").Value = "2.660.16"
$ws.Range("E17").Value = "  +3.53%  "

$ws.Range("D18").Value = "'12.59"
$ws.Range("E18").Value = "  +4.97%  "

$ws.Range("D19").Value = "'4.86"
$ws.Range("E19").Value = "  +5.62%  "

$ws.Range("D20").Value = "'359.71"
$ws.Range("E20").Value = "  +5.51%  "

$ws.Range("D21").Value = "'7.37"
$ws.Range("E21").Value = "  +9.67%  "

$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").Value = "'68.78"
$ws.Range("E23").Value = "  +4.28%  "

$ws.Range("D24").Value = "'9.45"
$ws.Range("E24").Value = "  +4.50%  "

$ws.Range("E25").Value = "  -2.18%  "

$ws.Range("D26").Value = "'1.67"
$ws.Range("E26").Value = "  +3.08%  "

$ws.Range("D27").Value = "'0.165"
$ws.Range("E27").Value = "  +2.98%  "

$ws.Range("D28").Value = "'8.12"
$ws.Range("E28").Value = "  +2.19%  "

$ws.Range("D29").Value = "0.0₃0956"
$ws.Range("E29").Value = "  +15.36%  "

$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").Value = "'533.17"
$ws.Range("E31").Value = "  -3.75%  "

$ws.Range("D32").Value = "'2.12"
$ws.Range("E32").Value = "  +6.79%  "

$ws.Range("D33").Value = "'1.79"
$ws.Range("E33").Value = "  +4.28%  "

$ws.Range("E34").Value = "  +6.59%  "

$ws.Range("D35").Value = "'6.38"
$ws.Range("E35").Value = "  +8.65%  "

$ws.Range("D36").Value = "'0.427"
$ws.Range("E36").Value = "  +5.53%  "

$ws.Range("D37").Value = "'20.46"
$ws.Range("E37").Value = "  +6.64%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'2.03"
$ws.Range("E38").Value = "  +9.05%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'162.97"
$ws.Range("E39").Value = "  -1.61%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.25%  "

$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").Value = "'42.38"
$ws.Range("E42").Value = "  +8.03%  "

$ws.Range("D43").Value = "'165.82"
$ws.Range("E43").Value = "  +0.93%  "

$ws.Range("D44").Value = "'4.13"
$ws.Range("E44").Value = "  +5.24%  "

$ws.Range("D45").Value = "'0.0619"
$ws.Range("E45").Value = "  +7.80%  "

$ws.Range("D46").Value = "'23.13"
$ws.Range("E46").Value = "  +4.11%  "

$ws.Range("D47").Value = "'2.25"
$ws.Range("E47").Value = "  +12.02%  "

$ws.Range("D48").Value = "'0.655"
$ws.Range("E48").Value = "  +5.88%  "

$ws.Range("D49").Value = "'0.0260"
$ws.Range("E49").Value = "  +6.42%  "

$ws.Range("D50").Value = "'0.0983"
$ws.Range("E50").Value = "  +3.49%  "

$ws.Range("D51").Value = "'19.51"
$ws.Range("E51").Value = "  +4.28%  "
